{"js": "// Adds the explicit Latin font (ascii/hAnsi = \"Arial\") to the two\n// signature-line (\"____________________________________\") paragraphs\n// in the \"Constancia de Liberaci\u00f3n de Actividades Complementarias\"\n// template, while keeping the existing complex-script Arial (w:cs)\n// untouched. Matches the author's fix (\"Resuelto: exceso de creditos\").\n//\n// - The \"{coordinador}\" signature line: both the paragraph mark\n//   (pPr/rPr) AND its run get ascii/hAnsi Arial added.\n// - The \"{jefe}\" signature line: only its run gets ascii/hAnsi Arial\n//   added; its paragraph mark formatting is left as-is.\n\nconst SIGNATURE_LINE = \"____________________________________\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate every signature-line paragraph and remember the text of the\n// paragraph right above it (that's what names whose signature it is),\n// instead of hard-coding paragraph indices.\nconst signatureParagraphs = [];\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === SIGNATURE_LINE) {\n    const precedingText = i > 0 ? items[i - 1].text : \"\";\n    signatureParagraphs.push({ paragraph: items[i], precedingText: precedingText });\n  }\n}\n\n// 1) \"{coordinador}\" signature line -> update the paragraph's font so\n//    both the paragraph mark and the run pick up ascii/hAnsi Arial\n//    (w:cs=\"Arial\" is preserved automatically).\nconst coordinadorEntry = signatureParagraphs.filter(function (entry) {\n  return entry.precedingText.indexOf(\"coordinador\") !== -1;\n})[0];\nif (coordinadorEntry) {\n  coordinadorEntry.paragraph.font.name = \"Arial\";\n}\n\n// 2) \"{jefe}\" signature line -> update only the run's font (ascii +\n//    hAnsi / \"Other\"), leaving the paragraph mark's rPr untouched.\nconst jefeEntry = signatureParagraphs.filter(function (entry) {\n  return entry.precedingText.indexOf(\"jefe\") !== -1;\n})[0];\nif (jefeEntry) {\n  const runRange = jefeEntry.paragraph.getRange(\"Content\");\n  const searchResults = runRange.search(SIGNATURE_LINE, { matchCase: true });\n  searchResults.load(\"items\");\n  await context.sync();\n  const runOnly = searchResults.items[0];\n  runOnly.font.nameAscii = \"Arial\";\n  runOnly.font.nameOther = \"Arial\";\n}\n\nawait context.sync();\n", "ps1": "# Adds the explicit Latin font (ascii/hAnsi = \"Arial\") to the two\n# signature-line (\"____________________________________\") paragraphs\n# in the \"Constancia de Liberaci\u00f3n de Actividades Complementarias\"\n# template, while keeping the existing complex-script Arial (w:cs)\n# untouched. Matches the author's fix (\"Resuelto: exceso de creditos\").\n#\n# - The \"{coordinador}\" signature line: both the paragraph mark and\n#   its run get ascii/hAnsi Arial added (Paragraphs.Item(i).Range\n#   includes the paragraph mark, so Font.Name there updates both).\n# - The \"{jefe}\" signature line: only its run gets ascii/hAnsi Arial\n#   added; its paragraph mark formatting is left as-is (so we use a\n#   Find range limited to the run's text, not the whole paragraph).\n\n$d = $word.ActiveDocument\n$SIGNATURE_LINE = \"____________________________________\"\n\n# Locate every signature-line paragraph and remember the text right\n# above it (that's what names whose signature it is), instead of\n# hard-coding paragraph indices.\n$paraCount = $d.Paragraphs.Count\n$signatureParaIndexes = @()\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $paraText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($paraText -eq $SIGNATURE_LINE) {\n        $signatureParaIndexes += $i\n    }\n}\n\n$coordinadorIndex = $null\n$jefeIndex = $null\nforeach ($idx in $signatureParaIndexes) {\n    $precedingText = \"\"\n    if ($idx -gt 1) {\n        $precedingText = $d.Paragraphs.Item($idx - 1).Range.Text\n    }\n    if ($precedingText -match \"coordinador\") {\n        $coordinadorIndex = $idx\n    } elseif ($precedingText -match \"jefe\") {\n        $jefeIndex = $idx\n    }\n}\n\n# 1) \"{coordinador}\" signature line -> the paragraph's own Range\n#    includes the paragraph mark, so setting Font.Name here updates\n#    both the paragraph mark (pPr/rPr) and the run (w:cs stays put).\nif ($coordinadorIndex) {\n    $coordinadorPara = $d.Paragraphs.Item($coordinadorIndex)\n    $coordinadorPara.Range.Font.Name = \"Arial\"\n}\n\n# 2) \"{jefe}\" signature line -> update only the run's font (ascii +\n#    \"other\"/hAnsi), leaving the paragraph mark's rPr untouched. Using\n#    Find scoped to that one paragraph's range keeps the match (and\n#    therefore the formatted range) limited to the run's text only.\nif ($jefeIndex) {\n    $jefePara = $d.Paragraphs.Item($jefeIndex)\n    $searchRange = $jefePara.Range.Duplicate\n    $find = $searchRange.Find\n    $find.ClearFormatting()\n    $find.Text = $SIGNATURE_LINE\n    $find.Forward = $true\n    $find.Wrap = 0\n    if ($find.Execute()) {\n        $searchRange.Font.NameAscii = \"Arial\"\n        $searchRange.Font.NameOther = \"Arial\"\n    }\n}\n"}
